# Update the training-data sample sheet with newly computed
# (MOSS-based) similarity scores and append one more training row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A / B refreshed values (column C is unchanged for existing rows)
$ws.Cells.Item(2, 2).Value  = 0.311

$ws.Cells.Item(3, 1).Value  = 0.425
$ws.Cells.Item(3, 2).Value  = 0.423

$ws.Cells.Item(4, 1).Value  = 0.003
$ws.Cells.Item(4, 2).Value  = 0.117

$ws.Cells.Item(5, 1).Value  = 0.013
$ws.Cells.Item(5, 2).Value  = 0

$ws.Cells.Item(6, 2).Value  = 0

$ws.Cells.Item(7, 1).Value  = 0.228
$ws.Cells.Item(7, 2).Value  = 0.228

$ws.Cells.Item(8, 1).Value  = 0.002
$ws.Cells.Item(8, 2).Value  = 0

$ws.Cells.Item(9, 1).Value  = 0.347
$ws.Cells.Item(9, 2).Value  = 0.326

$ws.Cells.Item(10, 1).Value = 0.182
$ws.Cells.Item(10, 2).Value = 0.171

$ws.Cells.Item(11, 2).Value = 1

$ws.Cells.Item(12, 1).Value = 0.002
$ws.Cells.Item(12, 2).Value = 0.039

$ws.Cells.Item(13, 1).Value = 0.049
$ws.Cells.Item(13, 2).Value = 0.172

$ws.Cells.Item(14, 2).Value = 0

$ws.Cells.Item(15, 1).Value = 0.01
$ws.Cells.Item(15, 2).Value = 0.181

$ws.Cells.Item(16, 1).Value = 0.103
$ws.Cells.Item(16, 2).Value = 0.042

$ws.Cells.Item(17, 1).Value = 0.012
$ws.Cells.Item(17, 2).Value = 0.009

$ws.Cells.Item(18, 1).Value = 0.03
$ws.Cells.Item(18, 2).Value = 0.105

# New training row 19
$ws.Cells.Item(19, 1).Value = 0.017
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0.77
